# Update "Horarios" workbook with newly scraped schedule data.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("LP1912")
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws3 = $wb.Worksheets.Item("6203-6173")

$newTime = "01:54:55"

# ---- Sheet 1: LP1912 ----
$ws1.Range("A2").Value = "Última actualización: $newTime"
$ws1.Range("A3").Value = "Total filas: 3"

# Row 6 (existing row updated)
$ws1.Cells.Item(6, 1).Value = $newTime
$ws1.Cells.Item(6, 2).Value = "01:58"
$ws1.Cells.Item(6, 3).Value = "14_ABASTO"
$ws1.Cells.Item(6, 4).Value = 4
$ws1.Cells.Item(6, 5).Value = "LP1912"

# Row 7 (existing row updated)
$ws1.Cells.Item(7, 1).Value = $newTime
$ws1.Cells.Item(7, 2).Value = "03:16"
$ws1.Cells.Item(7, 3).Value = "215_ALUAR"
$ws1.Cells.Item(7, 4).Value = 82
$ws1.Cells.Item(7, 5).Value = "LP1912"

# Row 8 (new row)
$ws1.Cells.Item(8, 1).Value = $newTime
$ws1.Cells.Item(8, 2).Value = "03:48"
$ws1.Cells.Item(8, 3).Value = "14_ABASTO"
$ws1.Cells.Item(8, 4).Value = 114
$ws1.Cells.Item(8, 5).Value = "LP1912"

# ---- Sheet 2: LP1912-215 ----
$ws2.Range("A2").Value = "Última actualización: $newTime"

$ws2.Cells.Item(6, 1).Value = $newTime
$ws2.Cells.Item(6, 2).Value = "03:16"
$ws2.Cells.Item(6, 3).Value = "215_ALUAR"
$ws2.Cells.Item(6, 4).Value = 82
$ws2.Cells.Item(6, 5).Value = "LP1912"

# ---- Sheet 3: 6203-6173 ----
$ws3.Range("A2").Value = "Última actualización: $newTime"
